$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4269.88162055857
$ws.Range("F2").Value = 1.3783800354989

# Row 3
$ws.Range("C3").Value = 4230.7412557336
$ws.Range("D3").Value = 4584
$ws.Range("F3").Value = -0.252469105797218

# Row 4
$ws.Range("C4").Value = 4238.32660211596
$ws.Range("F4").Value = 110.063586524759
